$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 51.873844
$ws.Range("H2").Value = 155.621532
$ws.Range("I2").Value = 0.9152980296207101
$ws.Range("J2").Value = 0.9152980296207101
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.807635
$ws.Range("N2").Value = 14.422905
$ws.Range("O2").Value = 0.1639819574772189
$ws.Range("P2").Value = 0.1639819574772189
$ws.Range("Q2").Value = 249.39050799894
$ws.Range("R2").Value = 2244.51457199046
$ws.Range("S2").Value = 0.1500923625722455
$ws.Range("T2").Value = 0.1500923625722455

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 51.873844
$ws.Range("H3").Value = 155.621532
$ws.Range("I3").Value = 0.9152980296207101
$ws.Range("J3").Value = 0.9152980296207101
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 18.86426066666667
$ws.Range("N3").Value = 56.592782
$ws.Range("O3").Value = 0.643434534959602
$ws.Range("P3").Value = 0.6434345349596021
$ws.Range("Q3").Value = 978.5617149980027
$ws.Range("R3").Value = 8807.055434982025
$ws.Range("S3").Value = 0.5889343620384416
$ws.Range("T3").Value = 0.5889343620384417

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 51.873844
$ws.Range("H4").Value = 155.621532
$ws.Range("I4").Value = 0.9152980296207101
$ws.Range("J4").Value = 0.9152980296207101
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.646177333333333
$ws.Range("N4").Value = 16.938532
$ws.Range("O4").Value = 0.192583507563179
$ws.Range("P4").Value = 0.192583507563179
$ws.Range("Q4").Value = 292.8889221856693
$ws.Range("R4").Value = 2636.000299671024
$ws.Range("S4").Value = 0.1762713050100228
$ws.Range("T4").Value = 0.1762713050100228

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.0716
$ws.Range("H5").Value = 9.2148
$ws.Range("I5").Value = 0.05419743768715064
$ws.Range("J5").Value = 0.05419743768715064
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.807635
$ws.Range("N5").Value = 14.422905
$ws.Range("O5").Value = 0.1639819574772189
$ws.Range("P5").Value = 0.1639819574772189
$ws.Range("Q5").Value = 14.767131666
$ws.Range("R5").Value = 132.904184994
$ws.Range("S5").Value = 0.008887401922188558
$ws.Range("T5").Value = 0.008887401922188558

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.0716
$ws.Range("H6").Value = 9.2148
$ws.Range("I6").Value = 0.05419743768715064
$ws.Range("J6").Value = 0.05419743768715064
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 18.86426066666667
$ws.Range("N6").Value = 56.592782
$ws.Range("O6").Value = 0.643434534959602
$ws.Range("P6").Value = 0.6434345349596021
$ws.Range("Q6").Value = 57.94346306373333
$ws.Range("R6").Value = 521.4911675736
$ws.Range("S6").Value = 0.03487250311423378
$ws.Range("T6").Value = 0.03487250311423379

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.0716
$ws.Range("H7").Value = 9.2148
$ws.Range("I7").Value = 0.05419743768715064
$ws.Range("J7").Value = 0.05419743768715064
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.646177333333333
$ws.Range("N7").Value = 16.938532
$ws.Range("O7").Value = 0.192583507563179
$ws.Range("P7").Value = 0.192583507563179
$ws.Range("Q7").Value = 17.34279829706667
$ws.Range("R7").Value = 156.0851846736
$ws.Range("S7").Value = 0.0104375326507283
$ws.Range("T7").Value = 0.0104375326507283

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.728822
$ws.Range("H8").Value = 5.186466
$ws.Range("I8").Value = 0.03050453269213933
$ws.Range("J8").Value = 0.03050453269213932
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.807635
$ws.Range("N8").Value = 14.422905
$ws.Range("O8").Value = 0.1639819574772189
$ws.Range("P8").Value = 0.1639819574772189
$ws.Range("Q8").Value = 8.31154515597
$ws.Range("R8").Value = 74.80390640373001
$ws.Range("S8").Value = 0.005002192982784825
$ws.Range("T8").Value = 0.005002192982784825

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.728822
$ws.Range("H9").Value = 5.186466
$ws.Range("I9").Value = 0.03050453269213933
$ws.Range("J9").Value = 0.03050453269213932
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 18.86426066666667
$ws.Range("N9").Value = 56.592782
$ws.Range("O9").Value = 0.643434534959602
$ws.Range("P9").Value = 0.6434345349596021
$ws.Range("Q9").Value = 32.612948854268
$ws.Range("R9").Value = 293.516539688412
$ws.Range("S9").Value = 0.01962766980692664
$ws.Range("T9").Value = 0.01962766980692664

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.728822
$ws.Range("H10").Value = 5.186466
$ws.Range("I10").Value = 0.03050453269213933
$ws.Range("J10").Value = 0.03050453269213932
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.646177333333333
$ws.Range("N10").Value = 16.938532
$ws.Range("O10").Value = 0.192583507563179
$ws.Range("P10").Value = 0.192583507563179
$ws.Range("Q10").Value = 9.761235589768001
$ws.Range("R10").Value = 87.851120307912
$ws.Range("S10").Value = 0.005874669902427854
$ws.Range("T10").Value = 0.005874669902427853
